$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 19 (keep header row 1 and the first two data rows)
$ws.Range("A4:B19").EntireRow.Delete()

# Update the remaining two data rows with new bin labels/values
$ws.Range("A2").Value = "10-15"
$ws.Range("B2").Value = 4.133277868047619

$ws.Range("A3").Value = "5-10"
$ws.Range("B3").Value = 3.903367674083333
